$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert new "LevelEditor" row at row 2 (category rows shift down by one) ---
$ws.Rows.Item(2).Insert()
# Copy formatting from the row that is now row 3 (a normal data row) so the new
# row matches the rest of the table instead of inheriting the header's style.
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = "LevelEditor"
$ws.Range("B2").Value = "Level Editor"
$ws.Range("C2").Value = "Editor de Niveles"

# --- 2) Insert four new pause-menu rows right before the "pause.ExitPopup.Title" row ---
# After step 1, the tooltip row is now row 49 and "pause.ExitPopup.Title" is row 50.
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(50).Insert()
}

$ws.Range("A50").Value = "pause.PlayLevel"
$ws.Range("B50").Value = "Play Level"
$ws.Range("C50").Value = "Jugar Nivel"

$ws.Range("A51").Value = "pause.NoSpawnObject"
$ws.Range("B51").Value = "There isn't any player spawn obj in the scene."
$ws.Range("C51").Value = "No hay ningún spawn del jugador en la escena."

$ws.Range("A52").Value = "pause.SaveLevel"
$ws.Range("B52").Value = "Save Level"
$ws.Range("C52").Value = "Guardar Nivel"

$ws.Range("A53").Value = "pause.NoChanges"
$ws.Range("B53").Value = "There are no changes to save."
$ws.Range("C53").Value = "No hay cambios para guardar."

# --- 3) Update the (now shifted) ExitPopup.Content English text (row 55) ---
$ws.Range("B55").Value = "Warning, exiting now will delete any changes you have made since the last save. Are you sure you want to continue?"

# --- 4) Restore view state (zoom / scroll / selection) to match the saved workbook ---
$ws.Application.ActiveWindow.Zoom = 130
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("C60").Select()
